# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.449.40"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.95"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.20"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2791"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06503"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.875.22"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07439"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.28"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.077"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.23"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6421"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.425.57"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.99"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "233.58"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007525"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.111.51"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.151"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.090"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.329"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.91"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.919"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1026"
$ws.Range("E29").Value = "  +9.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.376"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.270"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.002"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04981"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.175"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7412"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01927"
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9239"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.053"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.03"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9957"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4194"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.582"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.226"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.82"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1228"
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.855"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.438"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.56"
$ws.Range("E51").Value = "  -1.86%  "
